$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expenses")

# Row 3 ("6582cd83..." expense) got re-synced with newer data from the
# source DB: new _id / userId, corrected amount, a real description, and a
# refreshed createdAt/updatedAt timestamp.
$ws.Range("A3").Value = '"658d37fcd3cb29a14b7efda2"'
$ws.Range("B3").Value = "658d37d9d3cb29a14b7efd99"
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = "this is descii"
$ws.Range("F3").Value = 45288.3718136574
$ws.Range("G3").Value = 45288.3718136574

# The second expense row (row 4, the "tea"/10 entry's duplicate id 6582d2d3...)
# no longer exists in the export, so drop the whole row.
$ws.Rows.Item(4).Delete()
